$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 242, pushing the existing rows 242:247
# (old "Clementina"/"Clemenuless" data dated 44357) down to 245:250.
$ws.Rows("242:244").Insert()

# Common values shared by every row in this block.
$company  = "Comercializadora del Agro de Limarí"
$region   = "Coquimbo"
$periodo  = 4
$tipo     = "Fruta"
$rubroId  = 100102
$rubro    = "Cítricos"
$espId    = 100102004
$especie  = "Mandarina"
$unidad   = "`$/bandeja 10 kilos"
$zona     = "Provincia de Limarí"
$muestras = 10

# New row 242: Murcott / Especial
$ws.Range("A242").Value = 2
$ws.Range("B242").Value = $company
$ws.Range("C242").Value = $region
$ws.Range("D242").Value = 44461
$ws.Range("E242").Value = $periodo
$ws.Range("F242").Value = $tipo
$ws.Range("G242").Value = $rubroId
$ws.Range("H242").Value = $rubro
$ws.Range("I242").Value = $espId
$ws.Range("J242").Value = $especie
$ws.Range("K242").Value = "Murcott"
$ws.Range("L242").Value = "Especial"
$ws.Range("M242").Value = 400
$ws.Range("N242").Value = 4500
$ws.Range("O242").Value = 5000
$ws.Range("P242").Value = 4750
$ws.Range("Q242").Value = $unidad
$ws.Range("R242").Value = $zona
$ws.Range("S242").Value = 475
$ws.Range("T242").Value = $muestras

# New row 243: Murcott / Primera
$ws.Range("A243").Value = 2
$ws.Range("B243").Value = $company
$ws.Range("C243").Value = $region
$ws.Range("D243").Value = 44461
$ws.Range("E243").Value = $periodo
$ws.Range("F243").Value = $tipo
$ws.Range("G243").Value = $rubroId
$ws.Range("H243").Value = $rubro
$ws.Range("I243").Value = $espId
$ws.Range("J243").Value = $especie
$ws.Range("K243").Value = "Murcott"
$ws.Range("L243").Value = "Primera"
$ws.Range("M243").Value = 600
$ws.Range("N243").Value = 3500
$ws.Range("O243").Value = 4000
$ws.Range("P243").Value = 3750
$ws.Range("Q243").Value = $unidad
$ws.Range("R243").Value = $zona
$ws.Range("S243").Value = 375
$ws.Range("T243").Value = $muestras

# New row 244: Murcott / Segunda
$ws.Range("A244").Value = 2
$ws.Range("B244").Value = $company
$ws.Range("C244").Value = $region
$ws.Range("D244").Value = 44461
$ws.Range("E244").Value = $periodo
$ws.Range("F244").Value = $tipo
$ws.Range("G244").Value = $rubroId
$ws.Range("H244").Value = $rubro
$ws.Range("I244").Value = $espId
$ws.Range("J244").Value = $especie
$ws.Range("K244").Value = "Murcott"
$ws.Range("L244").Value = "Segunda"
$ws.Range("M244").Value = 400
$ws.Range("N244").Value = 2500
$ws.Range("O244").Value = 3000
$ws.Range("P244").Value = 2750
$ws.Range("Q244").Value = $unidad
$ws.Range("R244").Value = $zona
$ws.Range("S244").Value = 275
$ws.Range("T244").Value = $muestras
